# Applies the cryptos.xlsx price/volume/ranking update described in the commit
# message: "Updated cryptos list on Sat Sep  7 14:28:24 UTC 2024 with GitHub Actions"
#
# The sheet stores Coin/Link/Price/Volume(1h) as plain text (inline strings), so
# every "Price" (column D) cell is forced back to Text format before the new
# value is written -- otherwise Excel would auto-convert numeric-looking values
# (e.g. "1.00", "0.550") into real numbers and silently drop the formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.420.99'
$ws.Range("E2").Value = '  -1.04%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.266.36'
$ws.Range("E3").Value = '  -2.75%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.24%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '495.39'
$ws.Range("E5").Value = '  -0.22%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.34'
$ws.Range("E6").Value = '  -1.56%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.529'
$ws.Range("E8").Value = '  -0.12%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.292.72'
$ws.Range("E9").Value = '  -1.99%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0949'
$ws.Range("E10").Value = '  +0.53%  '

# Row 11
$ws.Range("E11").Value = '  +2.08%  '

# Row 12
$ws.Range("E12").Value = '  +2.63%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.61'
$ws.Range("E13").Value = '  -2.50%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.683.61'
$ws.Range("E14").Value = '  -2.39%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.81'
$ws.Range("E15").Value = '  +1.96%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '54.412.37'
$ws.Range("E16").Value = '  -1.03%  '

# Row 17
$ws.Range("E17").Value = '  +0.08%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.297.69'
$ws.Range("E18").Value = '  -1.37%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.03'
$ws.Range("E19").Value = '  +3.41%  '

# Row 20
$ws.Range("E20").Value = '  +2.20%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '305.27'
$ws.Range("E21").Value = '  +0.03%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.45'
$ws.Range("E22").Value = '  +4.43%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.15%  '

# Row 24
$ws.Range("E24").Value = '  -1.80%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.05'
$ws.Range("E25").Value = '  -2.10%  '

# Row 26
$ws.Range("E26").Value = '  +0.39%  '

# Row 27
$ws.Range("E27").Value = '  +5.41%  '

# Row 28
$ws.Range("E28").Value = '  +1.64%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.397.39'
$ws.Range("E29").Value = '  -2.21%  '

# Row 30
$ws.Range("E30").Value = '  +0.10%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '171.08'
$ws.Range("E31").Value = '  +1.87%  '

# Row 32
$ws.Range("E32").Value = '  -1.32%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0686'
$ws.Range("E33").Value = '  -2.04%  '

# Row 34
$ws.Range("E34").Value = '  +3.33%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.05%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.993'
$ws.Range("E36").Value = '  -0.73%  '

# Row 37
$ws.Range("E37").Value = '  +1.04%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.60'
$ws.Range("E38").Value = '  +0.41%  '

# Row 39
$ws.Range("E39").Value = '  +3.21%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.865'
$ws.Range("E40").Value = '  +2.59%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.65'
$ws.Range("E41").Value = '  +0.83%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.54'
$ws.Range("E42").Value = '  -1.21%  '

# Row 43
$ws.Range("E43").Value = '  +1.34%  '

# Row 44
$ws.Range("E44").Value = '  +2.03%  '

# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '129.14'
$ws.Range("E45").Value = '  +4.38%  '

# Row 46
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.34'
$ws.Range("E46").Value = '  +0.34%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.81'
$ws.Range("E47").Value = '  +1.83%  '

# Row 48
$ws.Range("E48").Value = '  +1.06%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.550'
$ws.Range("E49").Value = '  +0.22%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '242.67'
$ws.Range("E50").Value = '  +1.55%  '

# Row 51
$ws.Range("E51").Value = '  +1.51%  '
